# Updates the "Means" and "Standard Deviations" summary tables to include
# the tri-proximity (1/3/5/10 mile) columns' refreshed figures.

$wb = $excel.ActiveWorkbook

# --- Sheet: Means ---
$ws1 = $wb.Worksheets.Item("Means")

$ws1.Range("F2").Value = 70
$ws1.Range("G2").Value = 74

$ws1.Range("D3").Value = 14
$ws1.Range("E3").Value = 26
$ws1.Range("F3").Value = 24
$ws1.Range("G3").Value = 19

$ws1.Range("D4").Value = 4
$ws1.Range("E4").Value = 5.3
$ws1.Range("F4").Value = 6.4
$ws1.Range("G4").Value = 7.1

$ws1.Range("D5").Value = 6.9
$ws1.Range("E5").Value = 6.7
$ws1.Range("F5").Value = 7.8
$ws1.Range("G5").Value = 6.8

$ws1.Range("D6").Value = 74
$ws1.Range("E6").Value = 62
$ws1.Range("F6").Value = 55
$ws1.Range("G6").Value = 63

$ws1.Range("E7").Value = 8
$ws1.Range("F7").Value = 9.3
$ws1.Range("G7").Value = 8

$ws1.Range("D8").Value = 5.7
$ws1.Range("E8").Value = 6.7
$ws1.Range("F8").Value = 8
$ws1.Range("G8").Value = 7

$ws1.Range("D9").Value = 44

$ws1.Range("D10").Value = 0.47
$ws1.Range("G10").Value = 0.39

# --- Sheet: Standard Deviations ---
$ws2 = $wb.Worksheets.Item("Standard Deviations")

$ws2.Range("D2").Value = 24
$ws2.Range("E2").Value = 35
$ws2.Range("F2").Value = 31

$ws2.Range("D3").Value = 24
$ws2.Range("E3").Value = 36
$ws2.Range("F3").Value = 31

$ws2.Range("D4").Value = 3.7
$ws2.Range("E4").Value = 7
$ws2.Range("F4").Value = 7.9
$ws2.Range("G4").Value = 8.7

$ws2.Range("D5").Value = 14

$ws2.Range("D6").Value = 27
$ws2.Range("E6").Value = 26
$ws2.Range("F6").Value = 24
$ws2.Range("G6").Value = 29

$ws2.Range("D7").Value = 5.8
$ws2.Range("E7").Value = 9.9

$ws2.Range("D8").Value = 8.9
$ws2.Range("E8").Value = 9.4
$ws2.Range("G8").Value = 9.8

$ws2.Range("D9").Value = 32
$ws2.Range("E9").Value = 27
$ws2.Range("F9").Value = 18
$ws2.Range("G9").Value = 14

$ws2.Range("D10").Value = 0.41
$ws2.Range("E10").Value = 0.23
$ws2.Range("F10").Value = 0.17
$ws2.Range("G10").Value = 0.13

$wb.Save()
